$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11, shifting existing rows 11-16 down to 12-17
$ws.Rows.Item(11).Insert()

# Copy the style (number format) of the date cell from the row below (now row 12) to the new row 11
$ws.Range("D12").Copy()
$ws.Range("D11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate new row 11 with the weekly data point
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44455
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = 100112010
$ws.Cells.Item(11, 7).Value = "Achicoria"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 52
$ws.Cells.Item(11, 11).Value = 5000
$ws.Cells.Item(11, 12).Value = 6000
$ws.Cells.Item(11, 13).Value = 5500
$ws.Cells.Item(11, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(11, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(11, 16).Value = 344
$ws.Cells.Item(11, 17).Value = 16
$ws.Cells.Item(11, 18).Value = "Hortaliza"

$wb.Save()
